# Updated cryptos list (price + 1h volume-change columns) per the commit diff.
# A handful of D-column price cells are numeric-looking text with a trailing
# zero (e.g. "1.00", "0.140", "326.90"); a plain .Value assignment would let
# Excel silently coerce the string into a real number and drop that trailing
# zero. For just those cells we briefly force Text format ('@') so the literal
# string sticks, then restore NumberFormat/Style back to the original "Normal"
# so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.484.49"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "3.471.98"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "594.06"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").Value = "181.07"
$ws.Range("E6").Value = "  +3.73%  "
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  +6.64%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "3.468.75"
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.140"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.32%  "
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("D12").Value = "0.429"
$ws.Range("E12").Value = "  +1.84%  "
$ws.Range("D13").Value = "4.073.25"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").Value = "31.91"
$ws.Range("E14").Value = "  +3.52%  "
$ws.Range("D15").Value = "0.133"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").Value = "67.446.03"
$ws.Range("E17").Value = "  +1.86%  "
$ws.Range("D18").Value = "3.472.45"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").Value = "6.19"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "14.07"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").Value = "394.01"
$ws.Range("E21").Value = "  +2.43%  "
$ws.Range("E22").Value = "  +2.07%  "
$ws.Range("E23").Value = "  +1.41%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "0.998"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").Value = "0.539"
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("D26").Value = "71.85"
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("D28").Value = "10.37"
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").Value = "6.13"
$ws.Range("E31").Value = "  +1.29%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("D34").Value = "23.54"
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("D35").Value = "7.33"
$ws.Range("E35").Value = "  +1.83%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("D38").Value = "161.86"
$ws.Range("E38").Value = "  -0.85%  "
$ws.Range("D39").Value = "0.891"
$ws.Range("E39").Value = "  +2.88%  "
$ws.Range("D40").Value = "2.85"
$ws.Range("E40").Value = "  +11.27%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "4.68"
$ws.Range("E42").Value = "  +2.38%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "6.76"
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").Value = "26.18"
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("D45").Value = "0.0716"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "26.29"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.743.66"
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").Value = "41.55"
$ws.Range("E48").Value = "  -1.25%  "
$ws.Range("D49").Value = "0.0298"
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "326.90"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.63%  "
$ws.Range("E51").Value = "  -1.66%  "
